$wb = $excel.ActiveWorkbook

$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

$zhMsg = "Handback file name: bujwj3q0.kel is different with handoff file name: 50c72679-b681-4576-bc42-946bf680f3ed.04bfeffd4148abbada153f7fe93387d4db140935.zh-cn."
$deMsg = "Handback file name: bujwj3q0.kel is different with handoff file name: 50c72679-b681-4576-bc42-946bf680f3ed.04bfeffd4148abbada153f7fe93387d4db140935.de-de."
$statusMsg = "Handback transform failed"

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B7").Value = $statusMsg
$overview.Range("C7").Value = $statusMsg
$zhcn.Range("C7").Value = $statusMsg
$dede.Range("C7").Value = $statusMsg

$zhcn.Range("L7").Value = $zhMsg
$dede.Range("L7").Value = $deMsg
